# Update TPM-derived NATMI metrics for Apoe-Lrp5 LR-pair sheet.
# New ligand/receptor expression values propagate into the specificity
# and edge-weight columns for every sending/target cluster combination.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# row -> hashtable of column letter -> new value
$updates = @{
    2  = @{ G=47.23036199999999;  H=141.691086;  I=0.3244251370417807; J=0.3244251370417807;
            M=15.01856033333333;  N=45.055681;   O=0.4908713633047416; P=0.4908713633047417;
            Q=709.3320412621739;  R=6383.988371359565; S=0.1592510093100265; T=0.1592510093100266 }
    3  = @{ G=47.23036199999999;  H=141.691086;  I=0.3244251370417807; J=0.3244251370417807;
            O=0.3099803572711625; P=0.3099803572711625;
            Q=447.9360908202479;  R=4031.424817382232;  S=0.100565419887957;  T=0.100565419887957 }
    4  = @{ G=47.23036199999999;  H=141.691086;  I=0.3244251370417807; J=0.3244251370417807;
            O=0.1991482794240958; P=0.1991482794240958;
            Q=287.778563016408;   R=2590.007067147671;  S=0.06460870784379713; T=0.06460870784379713 }
    5  = @{ I=0.4188548944674916; J=0.4188548944674916;
            M=15.01856033333333;  N=45.055681;   O=0.4908713633047416; P=0.4908713633047417;
            Q=915.7958596994155;  R=8242.162737294739;  S=0.2056038730741213; T=0.2056038730741213 }
    6  = @{ I=0.4188548944674916; J=0.4188548944674916;
            O=0.3099803572711625; P=0.3099803572711625;
            S=0.1298367898318081; T=0.1298367898318081 }
    7  = @{ I=0.4188548944674916; J=0.4188548944674916;
            O=0.1991482794240958; P=0.1991482794240958;
            S=0.08341423156156218; T=0.08341423156156218 }
    8  = @{ I=0.2567199684907278; J=0.2567199684907277;
            M=15.01856033333333;  N=45.055681;   O=0.4908713633047416; P=0.4908713633047417;
            Q=561.299598861963;   R=5051.696389757666;  S=0.1260164809205938; T=0.1260164809205938 }
    9  = @{ I=0.2567199684907278; J=0.2567199684907277;
            O=0.3099803572711625; P=0.3099803572711625;
            S=0.07957814755139736; T=0.07957814755139736 }
    10 = @{ I=0.2567199684907278; J=0.2567199684907277;
            O=0.1991482794240958; P=0.1991482794240958;
            S=0.05112534001873652; T=0.05112534001873652 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
